$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Normalize specific party-abbreviation cells to their fully-uppercase form.
# Order matches how the shared-string table grows (new unique strings are
# appended in first-use order as the old, now-unused strings are dropped).
$ws.Range("C76").Value = "TAVK"
$ws.Range("C55").Value = "MKAT"
$ws.Range("C65").Value = "RPSN"
$ws.Range("C67").Value = "RTRJP"
$ws.Range("C69").Value = "SASAP"
$ws.Range("C70").Value = "SASP"
$ws.Range("C75").Value = "SUNP"
$ws.Range("C77").Value = "THP"
$ws.Range("C43").Value = "GOKMK"
$ws.Range("C49").Value = "JASD"
$ws.Range("C52").Value = "MAKKK"

# Update the saved view state to match (scrolled position / selected cell).
$ws.Range("C17").Select()
